$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status text for the "Comparative CT-Based Bone Density and Screw Trajectories" row
$ws.Range("C6").Value = "Submitted after Minor Revision"

# Update the selected/active cell and scroll position on the sheet view
$ws.Range("C7").Select()
